# The workbook tracks 4 "dimension" mapping blocks (columns C, K, AQ, BI =
# tamano-empresa, balance, ratios, pyg) that were re-curated so they now
# look like the other "measure" mapping blocks (e.g. column AS = ano):
#   - row 3 (the "kind" row) changes from "dim" to "medida"
#   - row 4 (the "datatype" row) changes from "skos:Concept" to "xsd:int"
#   - row 5 (the external mapping-file pointer) is cleared entirely
# Column BI additionally gets its row-2 identifier re-prefixed from
# "iaest-dimension:pyg" to "iaest-measure:pyg".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C", "K", "AQ", "BI")

foreach ($col in $cols) {
    $ws.Range($col + "3").Value = "medida"
    $ws.Range($col + "4").Value = "xsd:int"
    $ws.Range($col + "5").Clear()
}

$ws.Range("BI2").Value = "iaest-measure:pyg"
